$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90 holds a numeric-looking code ("999") that must be stored as TEXT,
# matching how every other row (e.g. A88="777", A89="888") in this sheet is
# stored. Force text entry via the "@" number format, then clear the
# formatting back off so no stray style sticks to the cell.
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = "999"
$ws.Range("A90").ClearFormats()

$ws.Range("B90").Value = "Incompleto"
$ws.Range("C90").Value = "Outro"
$ws.Range("D90").Value = "Zerar"
